# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated data (per commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    3  = 3004
    7  = 1642
    9  = 81
    10 = 31
    11 = 1345
    13 = 480
    14 = 342
    15 = 8
    16 = 70
    20 = 104
    21 = 3107
    22 = 383
    23 = 105
    24 = 64
    27 = 89
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
